$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value2 = 'KFY'

$ws.Range("A16").Value2 = 'KMA'

$ws.Range("A17").Value2 = 'KFY'

$ws.Range("A18").Value2 = 'CNB'

$ws.Range("A19").Value2 = 'KFY'

$ws.Range("A20").Value2 = 'CNB'

$ws.Range("A21").Value2 = 'KI'

$ws.Range("A22").Value2 = 'KMA'

$ws.Range("A23").Value2 = 'KFY'

$ws.Range("A24").Value2 = 'CNB'

$ws.Range("A41").Value2 = 'KCH'
$ws.Range("D41").Value2 = 'Bakalářská práce z chemie'
$ws.Range("E41").Value2 = 'Bakalářská práce z chemie'

$ws.Range("A42").Value2 = 'KFY'
$ws.Range("D42").Value2 = 'Bakalářská práce z fyziky'
$ws.Range("E42").Value2 = 'Bakalářská práce z fyziky'

$ws.Range("A50").Value2 = 'KCH'
$ws.Range("D50").Value2 = 'Diplomová práce z chemie'
$ws.Range("E50").Value2 = 'Diplomová práce z chemie'

$ws.Range("A51").Value2 = 'KFY'
$ws.Range("D51").Value2 = 'Diplomová práce fyziky'
$ws.Range("E51").Value2 = 'Diplomová práce'

$ws.Range("A71").Value2 = 'CNB'

$ws.Range("A72").Value2 = 'KFY'

$ws.Range("A84").Value2 = 'KPRF'

$ws.Range("A85").Value2 = 'KFY'

$ws.Range("A86").Value2 = 'KMA'
$ws.Range("D86").Value2 = 'SZZ - matematika s didaktikou'
$ws.Range("E86").Value2 = 'Matematika s didaktikou'

$ws.Range("A87").Value2 = 'KGEO'
$ws.Range("D87").Value2 = 'SZZ - Geografie s didaktikou'
$ws.Range("E87").Value2 = 'Geografie s didaktikou pro ZŠ'

$ws.Range("A90").Value2 = 'KMA'
$ws.Range("D90").Value2 = 'SZZ - matematika s didaktikou'
$ws.Range("E90").Value2 = 'Matematika s didaktikou'

$ws.Range("A91").Value2 = 'KGEO'
$ws.Range("D91").Value2 = 'SZZ - Geografie s didaktikou'
$ws.Range("E91").Value2 = 'Geografie s didaktikou pro SŠ'

$ws.Range("A93").Value2 = 'KFY'
$ws.Range("D93").Value2 = 'Aplikovaná fyzika'
$ws.Range("E93").Value2 = 'Aplikovaná fyzika'

$ws.Range("A94").Value2 = 'KGEO'
$ws.Range("D94").Value2 = 'Aplikovaná geografie - Cestovní ruch'
$ws.Range("E94").Value2 = 'Aplikovaná geografie - Cestovní ruch'

$ws.Range("A96").Value2 = 'KI'
$ws.Range("D96").Value2 = 'Informační technologie'
$ws.Range("E96").Value2 = 'Informační technologie'

$ws.Range("A98").Value2 = 'KGEO'
$ws.Range("D98").Value2 = 'Aplikovaná geografie - Krajina a GIS'
$ws.Range("E98").Value2 = 'Aplikovaná geografie - Krajina a GIS'

$ws.Range("A104").Value2 = 'KFY'
$ws.Range("D104").Value2 = 'SZZ - Fyzika s didaktikou pro SŠ'
$ws.Range("E104").Value2 = 'Fyzika s didaktikou pro SŠ'

$ws.Range("A105").Value2 = 'KCH'
$ws.Range("D105").Value2 = 'Chemie a didaktika chemie pro SŠ'
$ws.Range("E105").Value2 = 'Chemie a didaktika chemie pro SŠ'

$ws.Range("A106").Value2 = 'KGEO'
$ws.Range("D106").Value2 = 'Geografie krajiny a GIS'
$ws.Range("E106").Value2 = 'Geografie krajiny a GIS'

$ws.Range("A107").Value2 = 'KBI'
$ws.Range("D107").Value2 = 'Biologie s didaktikou pro SŠ'
$ws.Range("E107").Value2 = 'Biologie s didaktikou pro SŠ'

$ws.Range("A110").Value2 = 'KFY'
$ws.Range("D110").Value2 = 'SZZ - numerická matematika'
$ws.Range("E110").Value2 = 'Numerická matematika'

$ws.Range("A112").Value2 = 'KGEO'
$ws.Range("D112").Value2 = 'SZZ - Reg. geografie a reg. rozvoj Česka'
$ws.Range("E112").Value2 = 'Regionální geografie a regionální rozvoj Česka'

$ws.Range("A113").Value2 = 'KMA'
$ws.Range("D113").Value2 = 'Matematika pro vzdělávání'
$ws.Range("E113").Value2 = 'Matematika pro vzdělávání'

$ws.Range("A115").Value2 = 'KCH'
$ws.Range("D115").Value2 = 'Chemie pro vzdělávání'
$ws.Range("E115").Value2 = 'Chemie pro vzdělávání'

$ws.Range("A116").Value2 = 'KFY'
$ws.Range("D116").Value2 = 'Fyzika pro vzdělávání'
$ws.Range("E116").Value2 = 'Fyzika pro vzdělávání'

$ws.Range("A117").Value2 = 'KGEO'
$ws.Range("D117").Value2 = 'Geografie pro vzdělávání'
$ws.Range("E117").Value2 = 'Geografie pro vzdělávání'

$ws.Range("A118").Value2 = 'KBI'
$ws.Range("D118").Value2 = 'SZZ - Biologie'
$ws.Range("E118").Value2 = 'Biologie'

$ws.Range("A119").Value2 = 'KFY'
$ws.Range("D119").Value2 = 'Fyzika'
$ws.Range("E119").Value2 = 'Fyzika'

$ws.Range("A120").Value2 = 'KCH'
$ws.Range("D120").Value2 = 'Analytická chemie'
$ws.Range("E120").Value2 = 'Analytická chemie'

$ws.Range("A121").Value2 = 'KFY'
$ws.Range("D121").Value2 = 'Fyzika'
$ws.Range("E121").Value2 = 'Fyzika'

$ws.Range("A122").Value2 = 'KGEO'
$ws.Range("D122").Value2 = 'SZZ- Regionální geografie Evropy a světa'
$ws.Range("E122").Value2 = 'Regionální geografie Evropy a světa'

$ws.Range("A127").Value2 = 'KCH'
$ws.Range("D127").Value2 = 'Chemie'
$ws.Range("E127").Value2 = 'Chemie'

$ws.Range("A128").Value2 = 'KFY'
$ws.Range("D128").Value2 = 'Fyzika'
$ws.Range("E128").Value2 = 'Fyzika'

$ws.Range("A129").Value2 = 'KGEO'
$ws.Range("D129").Value2 = 'SZZ - Geografie'
$ws.Range("E129").Value2 = 'Geografie'

$ws.Range("A130").Value2 = 'KMA'
$ws.Range("D130").Value2 = 'Matematika'
$ws.Range("E130").Value2 = 'Matematika'

$ws.Range("A131").Value2 = 'KBI'
$ws.Range("D131").Value2 = 'Biologie'
$ws.Range("E131").Value2 = 'Biologie'

$ws.Range("A132").Value2 = 'KI'
$ws.Range("D132").Value2 = 'SZZ - Informatika'
$ws.Range("E132").Value2 = 'Informatika'

$ws.Range("A133").Value2 = 'KMA'
$ws.Range("D133").Value2 = 'Matematická informatika'
$ws.Range("E133").Value2 = 'Matematická informatika'

$ws.Range("A134").Value2 = 'KFY'
$ws.Range("D134").Value2 = 'Elektronika a elektrotechnika'
$ws.Range("E134").Value2 = 'Elektronika a elektrotechnika'

$ws.Range("A135").Value2 = 'KGEO'
$ws.Range("D135").Value2 = 'SZZ - Obecná geografie'
$ws.Range("E135").Value2 = 'Obecná geografie'

$ws.Range("A155").Value2 = 'KFY'
$ws.Range("D155").Value2 = 'Nanotechnologie a nanomateriály'
$ws.Range("E155").Value2 = 'Nanotechnologie a nanomateriály, jejich příprava a charakterizace'

$ws.Range("A156").Value2 = 'KCH'
$ws.Range("D156").Value2 = 'Syntéza, technologie a analýza'
$ws.Range("E156").Value2 = 'Syntéza, technologie a analýza chemických látek a materiálů'

$ws.Range("A159").Value2 = 'KFY'
$ws.Range("D159").Value2 = 'Učitelství fyziky pro střední školy'
$ws.Range("E159").Value2 = 'Učitelství fyziky pro střední školy'

$ws.Range("A160").Value2 = 'KGEO'
$ws.Range("D160").Value2 = 'Učitelství geografie pro střední školy'
$ws.Range("E160").Value2 = 'Učitelství geografie pro střední školy'

$ws.Range("A161").Value2 = 'KCH'
$ws.Range("D161").Value2 = 'Učitelství chemie pro střední školy'
$ws.Range("E161").Value2 = 'Učitelství chemie pro střední školy'

$ws.Range("A162").Value2 = 'KBI'
$ws.Range("D162").Value2 = 'Učitelství biologie pro střední školy'
$ws.Range("E162").Value2 = 'Učitelství biologie pro střední školy'

$ws.Range("A163").Value2 = 'KI'
$ws.Range("D163").Value2 = 'Učitelství informatiky pro střední školy'
$ws.Range("E163").Value2 = 'Učitelství informatiky pro střední školy'

$ws.Range("A164").Value2 = 'KMA'
$ws.Range("D164").Value2 = 'Učitelství matematiky pro střední školy'
$ws.Range("E164").Value2 = 'Učitelství matematiky pro střední školy'
